$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (columns F:I) - set in the order in which the new
# shared-string entries first appear in the final workbook
$ws.Range("F1").Value = "Min Price"
$ws.Range("G1").Value = "Max Price"
$ws.Range("H1").Value = "Os"
$ws.Range("H2").Value = "Android"
$ws.Range("I1").Value = "num"

# Updated values on row 2 (Password + Product columns)
$ws.Range("B2").Value = "Manik1591@#1"
$ws.Range("C2").Value = "Mobiles"

# Remaining new numeric values
$ws.Range("F2").Value = 5000
$ws.Range("G2").Value = 10000
$ws.Range("I2").Value = 5

# B2 now also becomes a hyperlink (mirrors the existing A2 mailto-link
# convention) styled with the built-in Hyperlink style
$null = $ws.Hyperlinks.Add($ws.Range("B2"), "mailto:Manik1591@#1")
$ws.Range("B2").Style = "Hyperlink"

# Final selection ends up on C2
$null = $ws.Range("C2").Select()
